# Apply data updates to the "Single sponsor institution stats" sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row 12 (Helsinki University Hospital) with new figures ---
$ws.Range("B12").Value = 46
$ws.Range("C12").Value = 40
$ws.Range("D12").Value = 87
$ws.Range("E12").Value = 74.3
$ws.Range("F12").Value = 93.89999999999999

# --- Re-insert "Örebro University" and "Örebro University Hospital" in
#     alphabetical order (after "Odense University Hospital", row 26),
#     pushing the remaining rows (Oslo University Hospital ... Zealand
#     University Hospital) down by two rows, and drop the now-duplicated
#     rows that used to sit at the end of the table. ---
$ws.Rows.Item(27).Insert()
$ws.Rows.Item(27).Insert()

$ws.Range("A27").Value = "Örebro University"
$ws.Range("B27").Value = 17
$ws.Range("C27").Value = 12
$ws.Range("D27").Value = 70.59999999999999
$ws.Range("E27").Value = 46.9
$ws.Range("F27").Value = 86.7

$ws.Range("A28").Value = "Örebro University Hospital"
$ws.Range("B28").Value = 1
$ws.Range("C28").Value = 1
$ws.Range("D28").Value = 100
$ws.Range("E28").Value = 5.1
$ws.Range("F28").Value = 100

# The original "Örebro University" / "Örebro University Hospital" rows,
# which used to be rows 54-55, are now (after the two-row insert above)
# at rows 56-57; remove them since that data now lives at rows 27-28.
$ws.Rows.Item(56).Delete()
$ws.Rows.Item(56).Delete()
